# Applies the resume text edits described by the commit diff.
# Each change is performed as a literal (non-wildcard) Find/Replace over
# $d.Content so it works regardless of which run(s) the text currently
# lives in.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. "L.A.M.P. (Linux, Apache, MySQL, PHP(basics)) / ..." -> drop "(basics)" wording
Replace-Text "  PHP(" "  PHP"
Replace-Text "basics))" ")"

# 2. Skills list tweaks - drop trailing qualifiers
Replace-Text "MS Office Suite / Basic Java" "MS Office  / Basic Java"
Replace-Text "BASH Scripting / Oracle VirtualBox / M" "BASH  / Oracle VirtualBox / M"

# 3. Add "Backfills, and" to the Disconnects/Reconnects bullet
Replace-Text "Disconnects/Reconnects, Quality Assurance," "Disconnects/Reconnects, Backfills, and Quality Assurance,"

# 4. Rewrite the "Level 3 support" bullet
Replace-Text "Level 3 support :Troubleshooting, N" "Level 3 support : Troubleshoot n"
Replace-Text "etwork connectivity and User support" "etwork connectivity issues. "
Replace-Text " on MS Office Suite,  Windows 7 & 10." "Provide user-support in MS Office,  Windows 7 & 10."

# 5. Fix typo: Macrum -> Macrium
Replace-Text "Macrum Reflect" "Macrium Reflect"

# 6. Capitalize "Blueprints" and hyphenate "low-voltage"
Replace-Text "Reading and designing blueprints, low voltage wiring and Ohm’s law." "Reading and designing Blueprints, low-voltage wiring and Ohm’s law."

# 7. National Grid job title / date formatting
Replace-Text "National Grid (Help Desk) (" "National Grid (Computer Help Desk"
Replace-Text "January " ", January "

# 8. Space out the en dash between years
Replace-Text "2008- " "2008 - "

# 9. "Back up" -> "Backup"
Replace-Text "Provided disaster recovery for servers (Back up and Safeguard data " "Provided disaster recovery for servers (Backup and Safeguard data "

# 10. Replace security bullet wording
Replace-Text "Recognized, Identified and removed all current security threats; Troubleshoot network issues " "Vulnerability and threat assessment; Risk Management; Troubleshoot network issues "

# 11. Fix product name: Active@KilDisk -> Active@ KillDisk
Replace-Text "Used “Active@KilDisk” software to ensure the safety of sensitive data not being sacrificed " "Used “Active@ KillDisk” software to ensure the safety of sensitive data not being sacrificed "

# 12. Rework the data-migration bullet
Replace-Text "Provided data-migration to/from the MS Windows to Linux, and the Mac OS X Platform" "Data-migration to/from the Microsoft Windows to Linux, and the Mac OS X"

# 13. Split the CompTIA certifications line into two, each with a
#     certificate number appended after a tab stop.
Replace-Text "CompTIA Linux+ Certified, CompTIA A+ Certified" "CompTIA Linux+ Certified`t(COMP001007308786)`r CompTIA A+ Certified`t(COMP001007308786)"
